# CIERRE 26- AGO 23
# Update the "ARQUITECTO" vale and the "VALES DE INSENTIVOS" vale, and
# move the active selection on sheet 1 as recorded by the author.

$wb = $excel.ActiveWorkbook

# --- Sheet "ARQUITECTO        " -------------------------------------------
$wsArq = $wb.Worksheets.Item("ARQUITECTO        ")

# Fix spacing in the amount-in-words cell.
$wsArq.Range("A2").Value = "CIEN   MIL        PESOS 00/100 M.N."

# Fix the signer's name (spacing + spelling correction).
$wsArq.Range("C8").Value = "Arq. Rodolfo Higuera Velazco"

# Move the selection to F14, as left by the author when they saved.
$wsArq.Activate()
$wsArq.Range("F14").Select() | Out-Null

# --- Sheet "VALES DE INSENTIVOS" -------------------------------------------
$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# These two cells keep the same wording, nothing to change here.
# (kept for clarity / documentation of the vale contents)
# $wsVales.Range("A2").Value = "SEIS   MIL     PESOS 00/100 M.N."
# $wsVales.Range("A4").Value = "INCENTIVO DEL MES DE  JULIO     2023"
